$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This revision of the export drops four accounts that no longer belong in
# the "Saldo" list. Locate each one by its account number (column A) and
# remove the whole row. Rows are matched by value (not a fixed row number)
# so the deletes are correct regardless of exactly where each account
# currently sits; resolving the row numbers up front (before any deletes
# happen) and then removing them back-to-front keeps the remaining matches'
# row numbers valid as we go.
$accountsToRemove = @("008004835", "008054285", "004211911", "005993550")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

$rowsToDelete = @()
foreach ($account in $accountsToRemove) {
    $cell = $ws.Range("A1:A$lastRow").Find($account)
    if ($cell) {
        $rowsToDelete += $cell.Row
    }
}

$rowsToDelete = $rowsToDelete | Sort-Object -Descending
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}
